$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the data rows belonging to the first two workers
# (rows 16-27: worker "LUIS JAIRO BELTRAN SANJUAN" rows 16-21 and
#  worker "JUAN CARLOS NAVARRO SANCHEZ" rows 22-27). The remaining
# worker "JESUS MANUEL NAVARRO SEPULVEDA" rows (28-32) shift up to 16-20.
$ws.Range("16:27").Delete()

# The period column (E) for the remaining worker was listed newest-first
# (1909..1905); reverse it so it reads oldest-first (1905..1909).
$ws.Range("E16").Value = "1905"
$ws.Range("E17").Value = "1906"
$ws.Range("E18").Value = "1907"
$ws.Range("E19").Value = "1908"
$ws.Range("E20").Value = "1909"

# Update the summary figures above the table.
$ws.Range("E11").Value = 165625
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 5
